$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chapter 4")
$ws.Activate()

# Column B mirrors column A for rows 1-15, with row 11 intentionally
# showing the "correct" answer differing from the user's answer in A11.
for ($r = 1; $r -le 15; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $b = $ws.Cells.Item($r, 2)
    if ($r -eq 11) {
        $b.Value = "A"
    } else {
        $b.Value = $a.Value2
    }
}

# Highlight row 11 (A11:B11) with the red font color used elsewhere
# in this workbook to flag an incorrect answer.
$ws.Range("A11:B11").Font.Color = 255

$ws.Range("E13").Select()
